# https://jira.hl7.org/browse/FHIR-49216 - Corrected prefix application in modular diagram
#
# 1) Rename the four leaf items in the "modular" diagram on slide 1 so the
#    prefix (A/B) comes first: ItemA-1 -> AItem1, ItemA-2 -> AItem2,
#    ItemB-1 -> BItem1, ItemB-2 -> BItem2. Only the first paragraph (the
#    item's own label) of each shape changes; the bullet paragraphs below
#    it are left untouched.
# 2) Refresh the cached "today" text of the datetimeFigureOut date field
#    that appears on every slide layout and on the slide master.

$p = $ppt.ActivePresentation

function Update-ItemLabel($shape, [string]$newText) {
    # Update only the first run of the first paragraph so sibling
    # paragraphs (bullets) keep their own runs/formatting untouched.
    $para = $shape.TextFrame.TextRange.Paragraphs(1, 1)
    $run = $para.Runs(1, 1)
    $run.Text = $newText
}

function Update-DateField($shapes) {
    for ($shpIdx = 1; $shpIdx -le $shapes.Count; $shpIdx++) {
        $sh = $shapes.Item($shpIdx)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                $sh.TextFrame.TextRange.Runs(1, 1).Text = "8/6/2025"
            }
        }
    }
}

# --- 1. Rename the four quadrant items -------------------------------------------------
$s = $p.Slides.Item(1)
$root = $s.Shapes.Item(1)

$renames = @{
    "Rectangle 17" = "AItem1"   # was ItemA-1
    "Rectangle 18" = "AItem2"   # was ItemA-2
    "Rectangle 19" = "BItem1"   # was ItemB-1
    "Rectangle 20" = "BItem2"   # was ItemB-2
}

for ($giIdx = 1; $giIdx -le $root.GroupItems.Count; $giIdx++) {
    $item = $root.GroupItems.Item($giIdx)
    if ($renames.ContainsKey($item.Name)) {
        Update-ItemLabel $item $renames[$item.Name]
    }
}

# --- 2. Refresh the cached date field text on every layout + the master ----------------
$master = $p.SlideMaster
Update-DateField $master.Shapes

for ($lytIdx = 1; $lytIdx -le $master.CustomLayouts.Count; $lytIdx++) {
    $layout = $master.CustomLayouts.Item($lytIdx)
    Update-DateField $layout.Shapes
}
